$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# In each 4-row year group (quarters A/B/C/D), the "B quarter" row and the
# "C quarter" row had their contents (columns A:E) transposed; the "A
# quarter" and "D quarter" rows were left in place. Use copy/paste so that
# cell formatting (including the empty placeholder cell in column C) is
# preserved exactly, rather than only round-tripping through .Value.
$holdRow = 200
for ($row = 3; $row -le 67; $row += 4) {
    $firstRow  = $row
    $secondRow = $row + 1

    $firstRange  = $ws.Range("A$($firstRow):E$($firstRow)")
    $secondRange = $ws.Range("A$($secondRow):E$($secondRow)")
    $holdRange   = $ws.Range("A$($holdRow):E$($holdRow)")

    $firstRange.Copy()
    $holdRange.PasteSpecial()

    $secondRange.Copy()
    $firstRange.PasteSpecial()

    $holdRange.Copy()
    $secondRange.PasteSpecial()

    $holdRange.Clear()
}

# The "产销率" (F) and "销售量" (G) columns duplicated data already present
# in columns B and E, so they were removed.
$ws.Range("F1:G69").EntireColumn.Delete()
